$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.1904
$ws.Range("C8").Value = -11.7247
$ws.Range("B12").Value = 5.6935
$ws.Range("C12").Value = -14.6018
$ws.Range("C14").Value = -11.6759
$ws.Range("C22").Value = -11.1607
